# Finish sankey on funds flow diagram
# Rebuild the data table on Sheet1 to reflect the final (post-edit) state:
#  - insert a new "sequestration -> wildlife restoration" flow row
#  - rename the "wildlife fund" target to "wildlife"
#  - remove the per-state (AL / TX / AK) breakdown rows that are no longer used

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Final target table (including header) for A1:C12
$data = @(
    @("source", "target", "value"),
    @("pistols", "wildlife restoration", 190),
    @("firearms", "wildlife restoration", 190),
    @("ammunition", "wildlife restoration", 185),
    @("bows-arrows", "wildlife restoration", 47),
    @("sequestration", "wildlife restoration", 5.6),
    @("wildlife restoration", "usfw", 11.6),
    @("wildlife restoration", "grants", 3),
    @("wildlife restoration", "apportionments", 599),
    @("apportionments", "wildlife", 474),
    @("apportionments", "hunter ed", 120),
    @("apportionments", "enhanced", 8)
)

# Clear out the old range (which extended to row 20) before writing the new,
# shorter table so no stale cells are left behind.
$ws.Range("A1:C20").Clear()

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowNum = $i + 1
    $ws.Cells.Item($rowNum, 1).Value = $data[$i][0]
    $ws.Cells.Item($rowNum, 2).Value = $data[$i][1]
    $ws.Cells.Item($rowNum, 3).Value = $data[$i][2]
}

# Update the selected cell to match the saved view state
$ws.Range("C6").Select()

$wb.Save()
